$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" columns -> "_FV2410", "_new" columns -> "_FV2504"
$fields = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($fields[$i] + "_FV2410")
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($fields[$i] + "_FV2504")
}

# Turn the data range into a table (adds autofilter + table part)
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (pane split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
